$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh LR-pair TPM statistics (Ccl2-Ackr4) with new expression values,
# including four additional sending/target cluster combinations (rows 17-21).
$arr = New-Object "object[,]" 20,20

$arr[0,0] = "ECs"
$arr[0,1] = "Ccl2"
$arr[0,2] = "Ackr4"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 10.93086033333333
$arr[0,7] = 32.792581
$arr[0,8] = 0.02611891973042464
$arr[0,9] = 0.02622656909968252
$arr[0,10] = 1
$arr[0,11] = 0.3333333333333333
$arr[0,12] = 0.011782
$arr[0,13] = 0.035346
$arr[0,14] = 0.0185453160301082
$arr[0,15] = 0.0190618769471875
$arr[0,16] = 0.1287873964473333
$arr[0,17] = 1.159086568026
$arr[0,18] = 0.0004843836207657534
$arr[0,19] = 0.0004999276329250583

$arr[1,0] = "ECs"
$arr[1,1] = "Ccl2"
$arr[1,2] = "Ackr4"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 10.93086033333333
$arr[1,7] = 32.792581
$arr[1,8] = 0.02611891973042464
$arr[1,9] = 0.02622656909968252
$arr[1,10] = 2
$arr[1,11] = 0.6666666666666666
$arr[1,12] = 0.450137
$arr[1,13] = 1.350411
$arr[1,14] = 0.7085327552066554
$arr[1,15] = 0.7282682145116399
$arr[1,16] = 4.920384677865667
$arr[1,17] = 44.283462100791
$arr[1,18] = 0.01850611015961924
$arr[1,19] = 0.01909997665099194

$arr[2,0] = "ECs"
$arr[2,1] = "Ccl2"
$arr[2,2] = "Ackr4"
$arr[2,3] = "MuSCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 10.93086033333333
$arr[2,7] = 32.792581
$arr[2,8] = 0.02611891973042464
$arr[2,9] = 0.02622656909968252
$arr[2,10] = 1
$arr[2,11] = 0.5
$arr[2,12] = 0.051649
$arr[2,13] = 0.103298
$arr[2,14] = 0.08129749003896268
$arr[2,15] = 0.05570796596193557
$arr[2,16] = 0.5645680053563333
$arr[2,17] = 3.387408032138
$arr[2,18] = 0.002123402616612663
$arr[2,19] = 0.001461028818703465

$arr[3,0] = "ECs"
$arr[3,1] = "Ccl2"
$arr[3,2] = "Ackr4"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 10.93086033333333
$arr[3,7] = 32.792581
$arr[3,8] = 0.02611891973042464
$arr[3,9] = 0.02622656909968252
$arr[3,10] = 1
$arr[3,11] = 0.3333333333333333
$arr[3,12] = 0.1217406666666667
$arr[3,13] = 0.365222
$arr[3,14] = 0.1916244387242736
$arr[3,15] = 0.1969619425792371
$arr[3,16] = 1.330730224220222
$arr[3,17] = 11.976572017982
$arr[3,18] = 0.005005023333426978
$arr[3,19] = 0.005165635997062062

$arr[4,0] = "FAPs"
$arr[4,1] = "Ccl2"
$arr[4,2] = "Ackr4"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 143.2163033333333
$arr[4,7] = 429.64891
$arr[4,8] = 0.3422104954945279
$arr[4,9] = 0.3436209192170106
$arr[4,10] = 1
$arr[4,11] = 0.3333333333333333
$arr[4,12] = 0.011782
$arr[4,13] = 0.035346
$arr[4,14] = 0.0185453160301082
$arr[4,15] = 0.0190618769471875
$arr[4,16] = 1.687374485873333
$arr[4,17] = 15.18637037286
$arr[4,18] = 0.006346401787765938
$arr[4,19] = 0.006550059678594114

$arr[5,0] = "FAPs"
$arr[5,1] = "Ccl2"
$arr[5,2] = "Ackr4"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 143.2163033333333
$arr[5,7] = 429.64891
$arr[5,8] = 0.3422104954945279
$arr[5,9] = 0.3436209192170106
$arr[5,10] = 2
$arr[5,11] = 0.6666666666666666
$arr[5,12] = 0.450137
$arr[5,13] = 1.350411
$arr[5,14] = 0.7085327552066554
$arr[5,15] = 0.7282682145116399
$arr[5,16] = 64.46695713355668
$arr[5,17] = 580.2026142020101
$arr[5,18] = 0.2424673452333726
$arr[5,19] = 0.2502481933070207

$arr[6,0] = "FAPs"
$arr[6,1] = "Ccl2"
$arr[6,2] = "Ackr4"
$arr[6,3] = "MuSCs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 143.2163033333333
$arr[6,7] = 429.64891
$arr[6,8] = 0.3422104954945279
$arr[6,9] = 0.3436209192170106
$arr[6,10] = 1
$arr[6,11] = 0.5
$arr[6,12] = 0.051649
$arr[6,13] = 0.103298
$arr[6,14] = 0.08129749003896268
$arr[6,15] = 0.05570796596193557
$arr[6,16] = 7.396978850863334
$arr[6,17] = 44.38187310518
$arr[6,18] = 0.02782085434869487
$arr[6,19] = 0.01914242247155024

$arr[7,0] = "FAPs"
$arr[7,1] = "Ccl2"
$arr[7,2] = "Ackr4"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 143.2163033333333
$arr[7,7] = 429.64891
$arr[7,8] = 0.3422104954945279
$arr[7,9] = 0.3436209192170106
$arr[7,10] = 1
$arr[7,11] = 0.3333333333333333
$arr[7,12] = 0.1217406666666667
$arr[7,13] = 0.365222
$arr[7,14] = 0.1916244387242736
$arr[7,15] = 0.1969619425792371
$arr[7,16] = 17.43524824533556
$arr[7,17] = 156.91723420802
$arr[7,18] = 0.06557589412469449
$arr[7,19] = 0.06768024375984551

$arr[8,0] = "Inflammatory-Mac"
$arr[8,1] = "Ccl2"
$arr[8,2] = "Ackr4"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 157.1889546666667
$arr[8,7] = 471.566864
$arr[8,8] = 0.3755976715691904
$arr[8,9] = 0.3771457008466821
$arr[8,10] = 1
$arr[8,11] = 0.3333333333333333
$arr[8,12] = 0.011782
$arr[8,13] = 0.035346
$arr[8,14] = 0.0185453160301082
$arr[8,15] = 0.0190618769471875
$arr[8,16] = 1.852000263882667
$arr[8,17] = 16.668002374944
$arr[8,18] = 0.006965577519423422
$arr[8,19] = 0.007189104940700244

$arr[9,0] = "Inflammatory-Mac"
$arr[9,1] = "Ccl2"
$arr[9,2] = "Ackr4"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 157.1889546666667
$arr[9,7] = 471.566864
$arr[9,8] = 0.3755976715691904
$arr[9,9] = 0.3771457008466821
$arr[9,10] = 2
$arr[9,11] = 0.6666666666666666
$arr[9,12] = 0.450137
$arr[9,13] = 1.350411
$arr[9,14] = 0.7085327552066554
$arr[9,15] = 0.7282682145116399
$arr[9,16] = 70.75656448678933
$arr[9,17] = 636.8090803811041
$arr[9,18] = 0.2661232530861229
$arr[9,19] = 0.2746632261663542

$arr[10,0] = "Inflammatory-Mac"
$arr[10,1] = "Ccl2"
$arr[10,2] = "Ackr4"
$arr[10,3] = "MuSCs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 157.1889546666667
$arr[10,7] = 471.566864
$arr[10,8] = 0.3755976715691904
$arr[10,9] = 0.3771457008466821
$arr[10,10] = 1
$arr[10,11] = 0.5
$arr[10,12] = 0.051649
$arr[10,13] = 0.103298
$arr[10,14] = 0.08129749003896268
$arr[10,15] = 0.05570796596193557
$arr[10,16] = 8.118652319578667
$arr[10,17] = 48.711913917472
$arr[10,18] = 0.03053514796305383
$arr[10,19] = 0.0210100198654573

$arr[11,0] = "Inflammatory-Mac"
$arr[11,1] = "Ccl2"
$arr[11,2] = "Ackr4"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 157.1889546666667
$arr[11,7] = 471.566864
$arr[11,8] = 0.3755976715691904
$arr[11,9] = 0.3771457008466821
$arr[11,10] = 1
$arr[11,11] = 0.3333333333333333
$arr[11,12] = 0.1217406666666667
$arr[11,13] = 0.365222
$arr[11,14] = 0.1916244387242736
$arr[11,15] = 0.1969619425792371
$arr[11,16] = 19.13628813375644
$arr[11,17] = 172.226593203808
$arr[11,18] = 0.07197369300059017
$arr[11,19] = 0.07428334987417032

$arr[12,0] = "MuSCs"
$arr[12,1] = "Ccl2"
$arr[12,2] = "Ackr4"
$arr[12,3] = "ECs"
$arr[12,4] = 2
$arr[12,5] = 1
$arr[12,6] = 5.153359
$arr[12,7] = 10.306718
$arr[12,8] = 0.01231377640537609
$arr[12,9] = 0.008243018499152039
$arr[12,10] = 1
$arr[12,11] = 0.3333333333333333
$arr[12,12] = 0.011782
$arr[12,13] = 0.035346
$arr[12,14] = 0.0185453160301082
$arr[12,15] = 0.0190618769471875
$arr[12,16] = 0.060716875738
$arr[12,17] = 0.364301254428
$arr[12,18] = 0.0002283628749617893
$arr[12,19] = 0.0001571274043042264

$arr[13,0] = "MuSCs"
$arr[13,1] = "Ccl2"
$arr[13,2] = "Ackr4"
$arr[13,3] = "FAPs"
$arr[13,4] = 2
$arr[13,5] = 1
$arr[13,6] = 5.153359
$arr[13,7] = 10.306718
$arr[13,8] = 0.01231377640537609
$arr[13,9] = 0.008243018499152039
$arr[13,10] = 2
$arr[13,11] = 0.6666666666666666
$arr[13,12] = 0.450137
$arr[13,13] = 1.350411
$arr[13,14] = 0.7085327552066554
$arr[13,15] = 0.7282682145116399
$arr[13,16] = 2.319717560183
$arr[13,17] = 13.918305361098
$arr[13,18] = 0.008724713923499825
$arr[13,19] = 0.006003128364563873

$arr[14,0] = "MuSCs"
$arr[14,1] = "Ccl2"
$arr[14,2] = "Ackr4"
$arr[14,3] = "MuSCs"
$arr[14,4] = 2
$arr[14,5] = 1
$arr[14,6] = 5.153359
$arr[14,7] = 10.306718
$arr[14,8] = 0.01231377640537609
$arr[14,9] = 0.008243018499152039
$arr[14,10] = 1
$arr[14,11] = 0.5
$arr[14,12] = 0.051649
$arr[14,13] = 0.103298
$arr[14,14] = 0.08129749003896268
$arr[14,15] = 0.05570796596193557
$arr[14,16] = 0.266165838991
$arr[14,17] = 1.064663355964
$arr[14,18] = 0.001001079114658076
$arr[14,19] = 0.000459201793974367

$arr[15,0] = "MuSCs"
$arr[15,1] = "Ccl2"
$arr[15,2] = "Ackr4"
$arr[15,3] = "Resolving-Mac"
$arr[15,4] = 2
$arr[15,5] = 1
$arr[15,6] = 5.153359
$arr[15,7] = 10.306718
$arr[15,8] = 0.01231377640537609
$arr[15,9] = 0.008243018499152039
$arr[15,10] = 1
$arr[15,11] = 0.3333333333333333
$arr[15,12] = 0.1217406666666667
$arr[15,13] = 0.365222
$arr[15,14] = 0.1916244387242736
$arr[15,15] = 0.1969619425792371
$arr[15,16] = 0.6273733602326667
$arr[15,17] = 3.764240161396
$arr[15,18] = 0.002359620492256396
$arr[15,19] = 0.001623560936309573

$arr[16,0] = "Resolving-Mac"
$arr[16,1] = "Ccl2"
$arr[16,2] = "Ackr4"
$arr[16,3] = "ECs"
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 102.0140613333333
$arr[16,7] = 306.042184
$arr[16,8] = 0.2437591368004809
$arr[16,9] = 0.2447637923374727
$arr[16,10] = 1
$arr[16,11] = 0.3333333333333333
$arr[16,12] = 0.011782
$arr[16,13] = 0.035346
$arr[16,14] = 0.0185453160301082
$arr[16,15] = 0.0190618769471875
$arr[16,16] = 1.201929670629334
$arr[16,17] = 10.817367035664
$arr[16,18] = 0.004520590227191295
$arr[16,19] = 0.00466565729066386

$arr[17,0] = "Resolving-Mac"
$arr[17,1] = "Ccl2"
$arr[17,2] = "Ackr4"
$arr[17,3] = "FAPs"
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 102.0140613333333
$arr[17,7] = 306.042184
$arr[17,8] = 0.2437591368004809
$arr[17,9] = 0.2447637923374727
$arr[17,10] = 2
$arr[17,11] = 0.6666666666666666
$arr[17,12] = 0.450137
$arr[17,13] = 1.350411
$arr[17,14] = 0.7085327552066554
$arr[17,15] = 0.7282682145116399
$arr[17,16] = 45.92030352640267
$arr[17,17] = 413.282731737624
$arr[17,18] = 0.1727113328040407
$arr[17,19] = 0.178253690022709

$arr[18,0] = "Resolving-Mac"
$arr[18,1] = "Ccl2"
$arr[18,2] = "Ackr4"
$arr[18,3] = "MuSCs"
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 102.0140613333333
$arr[18,7] = 306.042184
$arr[18,8] = 0.2437591368004809
$arr[18,9] = 0.2447637923374727
$arr[18,10] = 1
$arr[18,11] = 0.5
$arr[18,12] = 0.051649
$arr[18,13] = 0.103298
$arr[18,14] = 0.08129749003896268
$arr[18,15] = 0.05570796596193557
$arr[18,16] = 5.268924253805334
$arr[18,17] = 31.613545522832
$arr[18,18] = 0.01981700599594323
$arr[18,19] = 0.0136352930122502

$arr[19,0] = "Resolving-Mac"
$arr[19,1] = "Ccl2"
$arr[19,2] = "Ackr4"
$arr[19,3] = "Resolving-Mac"
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 102.0140613333333
$arr[19,7] = 306.042184
$arr[19,8] = 0.2437591368004809
$arr[19,9] = 0.2447637923374727
$arr[19,10] = 1
$arr[19,11] = 0.3333333333333333
$arr[19,12] = 0.1217406666666667
$arr[19,13] = 0.365222
$arr[19,14] = 0.1916244387242736
$arr[19,15] = 0.1969619425792371
$arr[19,16] = 12.41925983609422
$arr[19,17] = 111.773338524848
$arr[19,18] = 0.04671020777330558
$arr[19,19] = 0.0482091520118496

$ws.Range("A2:T21").Value = $arr
